$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'294.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.86%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'40.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.37%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.56%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07380"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.26%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.550"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.06%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9273"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.97%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.1197"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.80%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1778"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.65%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.04391"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.50%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08756"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.02%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.24%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001248"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.98%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005846"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.90%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.361"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.31%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.301"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.18%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3270"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.82%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.842"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.21%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1379"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.95%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-2.08%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03920"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.88%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001268"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.97%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.003827"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.55%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001234"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-3.75%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003721"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.26%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02354"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.00%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05071"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.80%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'20.04%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.55%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1297"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.98%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007380"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.09%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007311"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.91%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.2929"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-7.54%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006109"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-6.34%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.04667"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-81.47%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004199"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.25%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.19%"
$ws.Range("E51").Style = "Normal"
